$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The crawl was re-run later the same day; every scraped row's timestamp
# column (O) moves from the morning run to the afternoon run.
$ws.Range("O2:O398").Value = "2023-01-09 12:56:41"

# A handful of products picked up updated rating/review counts and price
# metadata between the two crawl passes.

# Delifrance Frischback Baguettes Classique mini 6 Stueck: one more rating
$ws.Range("D36").Value = 43

# Prix Garantie Frischback Mehrkorn Broetchen 6 Stueck: rating amount and
# average rating value both shifted
$ws.Range("D100").Value = 5
$ws.Range("E100").Value = 4.5

# Parisette: one more rating
$ws.Range("D148").Value = 17

# Pasquier Pitch Schokolade 8 Stueck: now flagged as out of stock online
$ws.Range("M203").Value = "Pasquier Pitch Schokolade 8 Stück - Online kein Bestand 30% ab 2 Aktion 4.50 Schweizer Franken"

# Betty Bossi Spitzbube 2x 80g: now flagged as out of stock online
$ws.Range("M225").Value = "Betty Bossi Spitzbube 2x  80g - Online kein Bestand 4.40 Schweizer Franken"

# Pasquier Pancakes Choco 10 Stueck: now flagged as out of stock online
$ws.Range("M338").Value = "Pasquier Pancakes Choco 10 Stück - Online kein Bestand 4.30 Schweizer Franken"
